$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed.
$ws.Range("F2").Value = -4
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 3
$ws.Range("F14").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("F27").Value = -2
$ws.Range("F42").Value = -5
$ws.Range("F43").Value = -2
